$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.435.93"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.862.00"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4770"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3765"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07332"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9351"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07799"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.891.08"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008888"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "27.503.39"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.114"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.938"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.024"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.322"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7597"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.732"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02057"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5705"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.31%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.117"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05286"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.078"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.705"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4920"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.010"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.665"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9183"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.45%  "
